# Update column F (dSF) values to re-pulled data per commit
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    3  = -6
    4  = -5
    5  = -2
    8  = -6
    10 = -5
    13 = -8
    18 = -2
    20 = -10
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
